# Clean up footnote markers (e.g. " [1]", " [5]") and embedded line-breaks
# in the text labels throughout the workbook. Line breaks inside a cell's
# text are collapsed to a single space, and trailing "[n]" citation
# markers are stripped, leaving a trailing space where they used to be.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $startRow = $used.Row
    $startCol = $used.Column
    $numRows = $used.Rows.Count
    $numCols = $used.Columns.Count
    $endRow = $startRow + $numRows - 1
    $endCol = $startCol + $numCols - 1

    for ($r = $startRow; $r -le $endRow; $r++) {
        for ($c = $startCol; $c -le $endCol; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            $val = $cell.Value2

            if ($val -ne $null -and $val -is [string]) {
                $newVal = $val -replace ' ?\[\d\]', ' '
                $newVal = $newVal -replace "`r`n", ' '
                $newVal = $newVal -replace "`n", ' '

                if ($newVal -ne $val) {
                    $cell.Value = $newVal
                }
            }
        }
    }
}
